$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string text: the 'questions' Python list literal re-serialized
# as indented JSON (this is what the edited sharedStrings.xml <t> contains).
$newText = 'questions = [
    {
        "title": "You''re building a high-throughput API for a cryptocurrency trading platform. Time is extremely important for this platform because microseconds count when processing high-volume trade orders. For communicating with the API, you want to choose the fastest verb for read-only operations. What HTTP verb should you choose for retrieving trade orders with the API server?",
        "ques_type": 2,
        "options": [
            "GET",
            "UPDATE",
            "POST",
            "DELETE"
        ],
        "score": "GET"
    },
    {
        "title": "You work for a customer relationship management (CRM) company. The company\u2019s clients gain CRM access through a RESTful API. The CRM allows clients to add contact information for customers, prospects, and related persons (e.g., virtual assistants or marketing directors). You want to choose an appropriate API request path so clients can easily retrieve information for a single contact while also being flexible for future software changes. Which of the following API paths should you use?",
        "ques_type": 2,
        "options": [
            "/customers/{customer_id}",
            "/contacts/{contact_id}",
            "/contacts/{contact_type}/all",
            "/customers/all"
        ],
        "score": "/contacts/{contact_id}"
    },
    {
        "title": "You work for a large social media network and need to manage error handling for the API. You\u2019re trying to decide on an appropriate error code for authentication failures based on nonexistent users and incorrect passwords. You want to balance security against brute-force attacks by providing descriptive and true error codes. Which HTTP error code(s) should you use to keep the system secure and still report that an error occurred?",
        "ques_type": 2,
        "options": [
            "404 if the user doesn\u2019t exist, and 403 if the password is wrong.",
            "403 if the user doesn\u2019t exist, and 401 if the password is wrong.",
            "500 if the user doesn\u2019t exist or if the password is wrong.",
            "401 if the user doesn\u2019t exist or if the password is wrong."
        ],
        "score": "401 if the user doesn\u2019t exist or if the password is wrong."
    },
    {
        "title": "You need to document a field for entering a telephone number into a user''s two-factor authentication settings. You want developers to understand what to place in the field. What is the most useful way to document this field?",
        "ques_type": 2,
        "options": [
            "{\n\u201cphone\u201d: \u201c123-123-1234\u201d\n}\n",
            "{\n\u201cphone\u201d: {phone number}\n}\n",
            "{\n\u201cphone\u201d: \u201cString\u201d\n}\n",
            "{\n\u201cphone\u201d: \u201cPhone Number\u201d\n}\n"
        ],
        "score": "{\n\u201cphone\u201d: \u201c123-123-1234\u201d\n}"
    }
]'

# A1 previously held a plain numeric 0 with a bold/border/centered style;
# A2 held the big JSON/text blob as a shared string. The edit removes A1's
# value+style and moves the text blob up into A1, then drops row 2 entirely.
$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $newText
$ws.Range("A1").Style = "Normal"
